$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 8 (shifts existing rows 8..36 down to 9..37),
# carrying the formatting of the row that is being pushed down.
$ws.Rows.Item(8).Insert()

# Re-apply the row formatting (fill/alignment) that row 8 should have,
# matching its neighbours (copy format from the row right below, which
# now holds what used to be row 8's style).
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A8:B8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New dictionary entry: a1-district -> district
$ws.Range("A8").Value = "a1-district"
$ws.Range("B8").Value = "district"

# Leave the selection where the author left it after editing.
$ws.Range("C14").Select() | Out-Null
